{"js": "const body = context.document.body;\n\n// Append the new \"Next week's Scrum Master\" bullet after the last\n// paragraph in the document (\"What did you learn as an individual? \").\n// Word automatically carries over the paragraph style / list numbering\n// (ListParagraph, numId 1, ilvl 0) from the paragraph it is inserted\n// next to, matching the surrounding bulleted list.\nconst newParagraph = body.insertParagraph(\n  \"Next week\\u2019s Scrum Master: Shawn\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The last paragraph in the body is \"What did you learn as an individual? \".\n# Insert a new bulleted paragraph right after it for the new Scrum Master note.\n$lastParagraph = $d.Paragraphs.Last\n$newRange = $lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Next week\u2019s Scrum Master: Shawn\"\n"}
